{"js": "// The edit replaces two pieces of text inside the second content\n// paragraph of the document:\n//   \"passport #\"  ->  \"national provider identifier\"\n//   \"121234567\"   ->  \"1144513326\"\n// (The rest of the diff is just proofing markup - <w:proofErr> spell/\n// grammar-check markers - splitting existing runs without changing any\n// visible text, so no further text edits are required.)\n\nconst body = context.document.body;\n\nconst idResults = body.search(\"passport #\", { matchCase: true, matchWholeWord: false });\nidResults.load(\"items\");\n\nconst numberResults = body.search(\"121234567\", { matchCase: true, matchWholeWord: false });\nnumberResults.load(\"items\");\n\nawait context.sync();\n\nfor (let i = 0; i < idResults.items.length; i++) {\n  idResults.items[i].insertText(\"national provider identifier\", Word.InsertLocation.replace);\n}\n\nfor (let i = 0; i < numberResults.items.length; i++) {\n  numberResults.items[i].insertText(\"1144513326\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The edit replaces two pieces of text inside the second content\n# paragraph of the document:\n#   \"passport #\"  ->  \"national provider identifier\"\n#   \"121234567\"   ->  \"1144513326\"\n# (The rest of the diff is just proofing markup - <w:proofErr> spell/\n# grammar-check markers - splitting existing runs without changing any\n# visible text, so no further text edits are required.)\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.Execute(\"passport #\", $false, $false, $false, $false, $false, $true, 1, $false, \"national provider identifier\", 2)\n\n$find2 = $d.Content.Find\n$find2.Execute(\"121234567\", $false, $false, $false, $false, $false, $true, 1, $false, \"1144513326\", 2)\n"}
